# Fix ITI values: replace the 3,3,4.5,6 repeating pattern (rows 2-61)
# with an alternating 2,4 pattern (rows 2-61).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 61; $r++) {
    if ($r % 2 -eq 0) {
        $ws.Cells.Item($r, 1).Value = 2
    } else {
        $ws.Cells.Item($r, 1).Value = 4
    }
}

# The saved diff shows the stale "B9" cell selection being cleared (the
# sheetView no longer carries a <selection> pointing away from the default).
# Re-select A1 so the persisted view state reflects the default selection.
$ws.Range("A1").Select()
